$d = $word.ActiveDocument

# Paragraphs that get a new trailing run containing a single space,
# inserted just before the paragraph mark (i.e. as the last run of the
# paragraph, after whatever content -- oMath or a literal "." run --
# was already there).
$spaceTargets = @(17, 18, 19, 28, 29, 30, 34, 35, 36)

foreach ($i in $spaceTargets) {
    $p = $d.Paragraphs($i)
    $r = $p.Range
    $body = $d.Range($r.Start, $r.End - 1)
    $body.InsertAfter(" ")
}

# Remove the trailing "Reference" paragraph (bold "Reference", line break,
# "Example equations from OpenStax ..." text, and the hyperlink to
# OpenStax) -- it was the last paragraph in the document, right before
# the section break.
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.Delete()
